$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record (Pepino dulce, Cultivar IV Región, Primera) was
# reported and must be inserted as row 156, pushing every following record
# (old rows 156-254) down by one row (new rows 157-255).
$ws.Rows.Item(156).Insert()

$ws.Cells.Item(156, 1).Value = 10
$ws.Cells.Item(156, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(156, 3).Value = "La Araucanía"
$ws.Cells.Item(156, 4).Value = 44762
$ws.Cells.Item(156, 5).Value = 9
$ws.Cells.Item(156, 6).Value = 100112043
$ws.Cells.Item(156, 7).Value = "Pepino dulce"
$ws.Cells.Item(156, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(156, 9).Value = "Primera"
$ws.Cells.Item(156, 10).Value = 110
$ws.Cells.Item(156, 11).Value = 21000
$ws.Cells.Item(156, 12).Value = 21000
$ws.Cells.Item(156, 13).Value = 21000
$ws.Cells.Item(156, 14).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(156, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(156, 16).Value = 1167
$ws.Cells.Item(156, 17).Value = 18
$ws.Cells.Item(156, 18).Value = "Hortaliza"
